$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: cardholder name and card number
$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces text storage (quotePrefix) so the long digit
# string is not reinterpreted as a number, matching the original inline-string type.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 01.07.2024"

# Row 6: transaction
$ws.Range("B6").Value = "05.07."
$ws.Range("C6").Value = "06.07."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 926863"
$ws.Range("E6").Value = "85,85-"

# Row 7: transaction
$ws.Range("B7").Value = "08.07."
$ws.Range("C7").Value = "09.07."
$ws.Range("D7").Value = "KARTENZ./08.07 LIDL RO"
$ws.Range("E7").Value = "147,15-"

# Row 8: transaction
$ws.Range("B8").Value = "10.07."
$ws.Range("C8").Value = "11.07."
$ws.Range("D8").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E8").Value = "45,93-"

# Row 9: cleared out (previously a transaction row, now empty)
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 14.07.2024"
$ws.Range("E12").Value = "278,93-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.07.2024"
